$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "yes" entries to "no" for testcases 2 and 4 (rows 3 and 5)
$ws.Range("B3").Value = "no"
$ws.Range("B5").Value = "no"

# Add a new row of data: execution status "done" for a new entry in column C
$ws.Range("C6").Value = "done"

# Update the active selection to reflect the new cell
$ws.Range("D5").Select()
